$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the leftover style (s="1") on D74:D75 left behind by the old blank-row
# formatting before giving them real values.
$ws.Range("D74:D75").ClearFormats()

# --- Column A: Country ---
$ws.Range("A74:A103").Value = "Canada"

# --- Column B: Colony ---
$ws.Range("B74").Value = "Machias Seal Island, NB"
$ws.Range("B75").Value = "Machias Seal Island, NB"
$ws.Range("B76").Value = "Machias Seal Island, NB"
$ws.Range("B77").Value = "Machias Seal Island, NB"
$ws.Range("B78").Value = "Machias Seal Island, NB"
$ws.Range("B79").Value = "Machias Seal Island, NB"
$ws.Range("B80").Value = "Machias Seal Island, NB"
$ws.Range("B81").Value = "Machias Seal Island, NB"
$ws.Range("B82").Value = "Machias Seal Island, NB"
$ws.Range("B83").Value = "Machias Seal Island, NB"
$ws.Range("B84").Value = "Machias Seal Island, NB"
$ws.Range("B85").Value = "Machias Seal Island, NB"
$ws.Range("B86").Value = "Machias Seal Island, NB"
$ws.Range("B87").Value = "Machias Seal Island, NB"
$ws.Range("B88").Value = "Machias Seal Island, NB"
$ws.Range("B89").Value = "Machias Seal Island, NB"
$ws.Range("B90").Value = "Machias Seal Island, NB"
$ws.Range("B91").Value = "Machias Seal Island, NB"
$ws.Range("B92").Value = "Machias Seal Island, NB"
$ws.Range("B93").Value = "Machias Seal Island, NB"
$ws.Range("B94").Value = "Machias Seal Island, NB"
$ws.Range("B95").Value = "Machias Seal Island, NB"
$ws.Range("B96").Value = "Baccalieu Island, NF"
$ws.Range("B97").Value = "Baccalieu Island, NF"
$ws.Range("B98").Value = "Baccalieu Island, NF"
$ws.Range("B99").Value = "Baccalieu Island, NF"
$ws.Range("B100").Value = "North Bird Island, NF"
$ws.Range("B101").Value = "North Bird Island, NF"
$ws.Range("B102").Value = "Puffin Islands (Lab), LB"
$ws.Range("B103").Value = "Puffin Islands (Lab), LB"

# --- Column C: Year ---
$ws.Range("C74").Value = 1965
$ws.Range("C75").Value = 1971
$ws.Range("C76").Value = 1974
$ws.Range("C77").Value = 1975
$ws.Range("C78").Value = 1976
$ws.Range("C79").Value = 1977
$ws.Range("C80").Value = 1978
$ws.Range("C81").Value = 1978
$ws.Range("C82").Value = 1979
$ws.Range("C83").Value = 1979
$ws.Range("C84").Value = 1980
$ws.Range("C85").Value = 1981
$ws.Range("C86").Value = 1982
$ws.Range("C87").Value = 1983
$ws.Range("C88").Value = 1987
$ws.Range("C89").Value = 1998
$ws.Range("C90").Value = 2000
$ws.Range("C91").Value = 2003
$ws.Range("C92").Value = 2011
$ws.Range("C93").Value = 2015
$ws.Range("C94").Value = 2016
$ws.Range("C95").Value = 2019
$ws.Range("C96").Value = 1979
$ws.Range("C97").Value = 1984
$ws.Range("C98").Value = 1996
$ws.Range("C99").Value = 2005
$ws.Range("C100").Value = 1985
$ws.Range("C101").Value = 2019
$ws.Range("C102").Value = 1978
$ws.Range("C103").Value = 2002

# --- Column D: Mature individuals ---
$ws.Range("D74").Value = 3000
$ws.Range("D75").Value = 3000
$ws.Range("D76").Value = 1780
$ws.Range("D77").Value = 1500
$ws.Range("D78").Value = 4600
$ws.Range("D79").Value = 1200
$ws.Range("D80").Value = 1100
$ws.Range("D81").Value = 3000
$ws.Range("D82").Value = 1300
$ws.Range("D83").Value = 1626
$ws.Range("D84").Value = 1500
$ws.Range("D85").Value = 1600
$ws.Range("D86").Value = 1600
$ws.Range("D87").Value = 1500
$ws.Range("D88").Value = 1700
$ws.Range("D89").Value = 2000
$ws.Range("D90").Value = 16048
$ws.Range("D91").Value = 14668
$ws.Range("D92").Value = 15676
$ws.Range("D93").Value = 10350
$ws.Range("D94").Value = 9592
$ws.Range("D95").Value = 17270
$ws.Range("D96").Value = 58748
$ws.Range("D97").Value = 60000
$ws.Range("D98").Value = 90600
$ws.Range("D99").Value = 150000
$ws.Range("D100").Value = 2000
$ws.Range("D101").Value = 18782
$ws.Range("D102").Value = 3894
$ws.Range("D103").Value = 1238

# --- Column J: Reference ---
$ws.Range("J74:J103").Value = "S. Wilhelm, ?"

# --- Column K: Notes ---
$ws.Range("K74:K103").Value = "No SE or other variance provided"

# Rows 100-103 (North Bird Island / Puffin Islands) pick up the same black-font
# style already used for several other rows (e.g. C15) - reuse it by copying format.
$ws.Range("C15").Copy()
$ws.Range("B100:D101").PasteSpecial(-4122)
$ws.Range("B102:C103").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Update the frozen-pane scroll position / active selection to match where the
#     user ended up after entering the new data ---
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$excel.ActiveWindow.ScrollRow = 86
$ws.Range("I109").Select()
